# Update the cryptocurrency Price (D) and Volume(1h) (E) columns with
# freshly scraped values. D is written with a leading apostrophe so Excel
# keeps the price strings (several of which are plain-number-looking, e.g.
# "1.004") as literal text instead of auto-converting them to numbers; the
# style is then reset to "Normal" so the text-quote marker does not linger
# as a visible/stored number-format override.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "2" = @("28.893.48", "  -0.83%  ")
    "3" = @("1.880.43", "  -1.30%  ")
    "4" = @("1.004", "  -0.05%  ")
    "5" = @("324.54", "  -1.07%  ")
    "6" = @("1.003", "  -0.15%  ")
    "7" = @("0.4609", "  -1.20%  ")
    "8" = @("0.3869", "  -1.59%  ")
    "9" = @("0.07847", "  -1.90%  ")
    "10" = @("0.9852", "  -2.92%  ")
    "11" = @("21.79", "  -2.02%  ")
    "12" = @("1.879.57", "  -1.99%  ")
    "13" = @("6.997", "  -2.14%  ")
    "14" = @("5.650", "  -2.36%  ")
    "15" = @("0.06980", "  -0.06%  ")
    "16" = @("88.08", "  -1.94%  ")
    "17" = @("1.004", "  -0.01%  ")
    "18" = @("0.000009957", "  -1.97%  ")
    "19" = @("16.91", "  -2.51%  ")
    "20" = @("1.003", "  -0.22%  ")
    "21" = @("28.900.19", "  -0.76%  ")
    "22" = @("5.255", "  -2.23%  ")
    "23" = @("10.96", "  -1.60%  ")
    "24" = @("2.106", "  +1.88%  ")
    "25" = @("156.48", "  +0.57%  ")
    "26" = @("19.35", "  -2.21%  ")
    "27" = @("5.961", "  +1.20%  ")
    "28" = @("117.61", "  -2.60%  ")
    "29" = @("1.910", "  -4.56%  ")
    "30" = @("0.09359", "  -0.45%  ")
    "31" = @("0.9005", "  -4.47%  ")
    "32" = @("5.262", "  -2.12%  ")
    "33" = @("1.320", "  -2.86%  ")
    "34" = @("3.254", "  -0.26%  ")
    "35" = @("1.172", "  -0.57%  ")
    "36" = @("0.05740", "  -2.16%  ")
    "37" = @("0.02074", "  -1.60%  ")
    "38" = @("1.002", "  -0.24%  ")
    "39" = @("7.663", "  -5.84%  ")
    "40" = @("0.5654", "  -3.60%  ")
    "41" = @("0.1765", "  -3.02%  ")
    "42" = @("9.720", "  -3.16%  ")
    "43" = @("2.247", "  -1.62%  ")
    "44" = @("11.89", "  -0.86%  ")
    "45" = @("0.5342", "  -2.45%  ")
    "46" = @("0.07046", "  -2.44%  ")
    "47" = @("1.841", "  -2.24%  ")
    "48" = @("2.553", "  +2.16%  ")
    "49" = @("112.65", "  -0.82%  ")
    "50" = @("1.060", "  -5.32%  ")
    "51" = @("70.84", "  -0.83%  ")
}

foreach ($row in $updates.Keys) {
    $price = $updates[$row][0]
    $volume = $updates[$row][1]
    $ws.Range("D$row").Value = "'" + $price
    $ws.Range("D$row").Style = "Normal"
    $ws.Range("E$row").Value = $volume
}
